{"js": "const paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nconst p = paras.items[0];\nconst range = p.getRange();\nrange.load(\"text\");\nawait context.sync();\nconst text = range.text;\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">\n<w:body>\n<w:p>\n<w:r>\n<w:rPr>\n<w:rFonts w:ascii=\"Helvetica\" w:cs=\"Helvetica\" w:eastAsia=\"Helvetica\" w:hAnsi=\"Helvetica\"/>\n<w:sz w:val=\"24\"/>\n<w:szCs w:val=\"24\"/>\n<w:spacing w:val=\"20\"/>\n<w14:ligatures w14:val=\"historicalDiscretional\"/>\n</w:rPr>\n<w:t xml:space=\"preserve\">${text}</w:t>\n</w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\nrange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\nreturn \"done\";\n", "ps1": "$d = $word.ActiveDocument\n$xml = $d.WordOpenXML\n$newXml = $xml.Replace('<w:style w:type=\"paragraph\" w:default=\"1\" w:styleId=\"Normal\"><w:name w:val=\"Normal\"/></w:style>', '<w:style w:type=\"paragraph\" w:default=\"1\" w:styleId=\"Normal\"><w:name w:val=\"Normal\"/><w:rPr><w14:ligatures w14:val=\"historicalDiscretional\"/></w:rPr></w:style>')\n$d.WordOpenXML = $newXml\nWrite-Output \"set done\"\n"}
